$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date in column C (rows 2-44) from 45732 (2025-03-16)
# to 45733 (2025-03-17), keeping the existing number format/style intact.
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 3).Value = 45733
}
